# ks2 + ks3 + some works
# Restructure "Раздел 1" subsection tables: subsection0 grows from 3 to 4 rows,
# subsection1 shrinks from 2 to 1 row; section/sub totals + the TST section
# become generic template placeholders; numeric values collapse to simple
# placeholder figures (1's) with recomputed totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Structural edit: insert a 4th data row into subsection0 of "Раздел 1"
#    (after current row 6), reusing row 6's formatting, then drop the
#    2nd data row of subsection1 (old row 10, which lands on row 11 once the
#    insert above has shifted things down) so the overall row count is
#    unchanged (44 rows / dimension untouched).
# ---------------------------------------------------------------------------

$ws.Rows.Item(7).Insert(-4121)
$ws.Range("A6:J6").Copy()
$ws.Range("A7:J7").PasteSpecial(-4122)
$ws.Rows.Item(7).RowHeight = $ws.Rows.Item(6).RowHeight
$ws.Range("A7:B7").Merge()
$ws.Range("C7:F7").Merge()

$ws.Rows.Item(11).Delete()

# ---------------------------------------------------------------------------
# 2) Content edit: rewrite rows 4-22 with the new template text / numbers.
# ---------------------------------------------------------------------------

# --- Раздел 1 / subsection "Монтажные и пусконаладочные работы по разделу 1:" ---
# 4 generic data rows (А=1..4), all "название" / "шт" / 1 / 1 / 1
for ($i = 0; $i -lt 4; $i++) {
    $r = 4 + $i
    $ws.Cells.Item($r, 1).Value = $i + 1
    $ws.Cells.Item($r, 3).Value = "название"
    $ws.Cells.Item($r, 7).Value = "шт"
    $ws.Cells.Item($r, 8).Value = 1
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = 1
}

# subtotal row for subsection0
$ws.Cells.Item(8, 1).Value = "Итого"
$ws.Cells.Item(8, 10).Value = 4

# subsection1 header "Оборудование и материалы по разделу 1:" stays the same text,
# only its row moved down to 9 (no text change needed - already correct)

# 1 generic data row (subsection1)
$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 3).Value = "название"
$ws.Cells.Item(10, 7).Value = "шт"
$ws.Cells.Item(10, 8).Value = 1
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(10, 10).Value = 1

# subtotal row for subsection1
$ws.Cells.Item(11, 1).Value = "Итого"
$ws.Cells.Item(11, 10).Value = 1

# grand total for "Раздел 1"
$ws.Cells.Item(12, 1).Value = "Итого по разделу"
$ws.Cells.Item(12, 10).Value = 5

# --- Раздел TST ---
$ws.Cells.Item(14, 1).Value = "Подраздел:"

$ws.Cells.Item(15, 1).Value = 6
$ws.Cells.Item(15, 3).Value = "название"
$ws.Cells.Item(15, 7).Value = "шт"
$ws.Cells.Item(15, 8).Value = 1
$ws.Cells.Item(15, 9).Value = 1
$ws.Cells.Item(15, 10).Value = 1

$ws.Cells.Item(16, 1).Value = "Итого"
$ws.Cells.Item(16, 10).Value = 1

$ws.Cells.Item(17, 1).Value = "Подраздел:"

$ws.Cells.Item(18, 1).Value = 7
$ws.Cells.Item(18, 3).Value = "название"
$ws.Cells.Item(18, 7).Value = "шт"
$ws.Cells.Item(18, 8).Value = 1
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(18, 10).Value = 1

$ws.Cells.Item(19, 1).Value = "Итого:"
$ws.Cells.Item(19, 10).Value = 1

$ws.Cells.Item(20, 1).Value = "Итого по разделу"
$ws.Cells.Item(20, 10).Value = 2

# --- Footer totals ---
$ws.Cells.Item(21, 10).Value = 0.56
$ws.Cells.Item(22, 10).Value = 9.072000000000001

Write-Output "edit applied"
